$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextCell $ws "D2" "318.33"
Set-TextCell $ws "E2" "3.68%"
Set-TextCell $ws "G2" "16"
Set-TextCell $ws "D3" "39.83"
Set-TextCell $ws "E3" "2.98%"
Set-TextCell $ws "G3" "16"
Set-TextCell $ws "D4" "5.138"
Set-TextCell $ws "E4" "1.04%"
Set-TextCell $ws "G4" "16"
Set-TextCell $ws "D5" "0.08201"
Set-TextCell $ws "E5" "1.75%"
Set-TextCell $ws "G5" "16"
Set-TextCell $ws "D6" "2.006"
Set-TextCell $ws "E6" "3.92%"
Set-TextCell $ws "G6" "16"
Set-TextCell $ws "D7" "8.269"
Set-TextCell $ws "E7" "4.15%"
Set-TextCell $ws "G7" "16"
Set-TextCell $ws "D8" "4.281"
Set-TextCell $ws "E8" "2.15%"
Set-TextCell $ws "G8" "16"
Set-TextCell $ws "D9" "0.9336"
Set-TextCell $ws "E9" "0.38%"
Set-TextCell $ws "G9" "16"
Set-TextCell $ws "E10" "-3.07%"
Set-TextCell $ws "G10" "16"
Set-TextCell $ws "D11" "0.1968"
Set-TextCell $ws "E11" "1.82%"
Set-TextCell $ws "G11" "16"
Set-TextCell $ws "D12" "0.09068"
Set-TextCell $ws "E12" "1.35%"
Set-TextCell $ws "G12" "16"
Set-TextCell $ws "D13" "0.03553"
Set-TextCell $ws "E13" "1.67%"
Set-TextCell $ws "G13" "16"
Set-TextCell $ws "D14" "0.09803"
Set-TextCell $ws "E14" "0.12%"
Set-TextCell $ws "G14" "16"
Set-TextCell $ws "D15" "0.001392"
Set-TextCell $ws "E15" "-0.18%"
Set-TextCell $ws "G15" "16"
Set-TextCell $ws "D16" "0.006402"
Set-TextCell $ws "E16" "6.31%"
Set-TextCell $ws "G16" "16"
Set-TextCell $ws "D17" "3.666"
Set-TextCell $ws "E17" "-2.03%"
Set-TextCell $ws "G17" "16"
Set-TextCell $ws "E18" "-8.35%"
Set-TextCell $ws "G18" "16"
Set-TextCell $ws "D19" "0.3460"
Set-TextCell $ws "E19" "-0.11%"
Set-TextCell $ws "G19" "16"
Set-TextCell $ws "D20" "0.1275"
Set-TextCell $ws "E20" "-3.45%"
Set-TextCell $ws "G20" "16"
Set-TextCell $ws "D21" "4.905"
Set-TextCell $ws "E21" "2.64%"
Set-TextCell $ws "G21" "16"
Set-TextCell $ws "E22" "-2.01%"
Set-TextCell $ws "G22" "16"
Set-TextCell $ws "D23" "0.04327"
Set-TextCell $ws "E23" "-0.97%"
Set-TextCell $ws "G23" "16"
Set-TextCell $ws "E24" "-0.78%"
Set-TextCell $ws "G24" "16"
Set-TextCell $ws "D25" "0.004771"
Set-TextCell $ws "E25" "11.50%"
Set-TextCell $ws "G25" "16"
Set-TextCell $ws "D26" "0.0001301"
Set-TextCell $ws "E26" "0.11%"
Set-TextCell $ws "G26" "16"
Set-TextCell $ws "D27" "0.0003997"
Set-TextCell $ws "E27" "-10.13%"
Set-TextCell $ws "G27" "16"
Set-TextCell $ws "G28" "16"
Set-TextCell $ws "G29" "16"
Set-TextCell $ws "G30" "16"
Set-TextCell $ws "G31" "16"
Set-TextCell $ws "G32" "16"
Set-TextCell $ws "G33" "16"
Set-TextCell $ws "G34" "16"
Set-TextCell $ws "G35" "16"
Set-TextCell $ws "G36" "16"
Set-TextCell $ws "G37" "16"
Set-TextCell $ws "G38" "16"
Set-TextCell $ws "D39" "0.02217"
Set-TextCell $ws "E39" "7.17%"
Set-TextCell $ws "G39" "16"
Set-TextCell $ws "D40" "0.05241"
Set-TextCell $ws "E40" "3.68%"
Set-TextCell $ws "G40" "16"
Set-TextCell $ws "D41" "0.007516"
Set-TextCell $ws "E41" "1.02%"
Set-TextCell $ws "G41" "16"
Set-TextCell $ws "D42" "0.01023"
Set-TextCell $ws "E42" "1.32%"
Set-TextCell $ws "G42" "16"
Set-TextCell $ws "D43" "0.1378"
Set-TextCell $ws "E43" "1.96%"
Set-TextCell $ws "G43" "16"
Set-TextCell $ws "D44" "0.002152"
Set-TextCell $ws "E44" "0.58%"
Set-TextCell $ws "G44" "16"
Set-TextCell $ws "D45" "0.009879"
Set-TextCell $ws "E45" "8.97%"
Set-TextCell $ws "G45" "16"
Set-TextCell $ws "D46" "0.00006654"
Set-TextCell $ws "E46" "7.31%"
Set-TextCell $ws "G46" "16"
Set-TextCell $ws "D47" "0.00000000751"
Set-TextCell $ws "E47" "0.11%"
Set-TextCell $ws "G47" "16"
Set-TextCell $ws "D48" "0.002769"
Set-TextCell $ws "E48" "-1.03%"
Set-TextCell $ws "G48" "16"
Set-TextCell $ws "E49" "-24.95%"
Set-TextCell $ws "G49" "16"
Set-TextCell $ws "D50" "0.00002102"
Set-TextCell $ws "E50" "0.11%"
Set-TextCell $ws "G50" "16"
Set-TextCell $ws "D51" "0.0002002"
Set-TextCell $ws "E51" "0.11%"
Set-TextCell $ws "G51" "16"
